$d = $word.ActiveDocument

# Remove the hidden _GoBack bookmark from its current location (end of the
# "Bestellarray..." paragraph); it will be re-added at the end of the
# newly-inserted last paragraph, matching the target document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The last paragraph currently holds:
#   "Bestellarray wurde verändert: Name:Anzahl => Name:Anzahl:Preis"
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

# Insert a new paragraph after it containing the date.
$lastPara.Range.InsertParagraphAfter()
$datePara = $d.Paragraphs($d.Paragraphs.Count)
$datePara.Range.Text = "15.05.2017"

# Insert another new paragraph after that one with the new note text.
$datePara.Range.InsertParagraphAfter()
$notePara = $d.Paragraphs($d.Paragraphs.Count)
$notePara.Range.Text = "Bestellungen für einzelne Tische gespeichert. Liste muss aktualisiert werden"

# Re-add the _GoBack bookmark at the end of the new last paragraph (as a
# collapsed bookmark, matching the original's placement right after the text).
$d.Bookmarks.Add("_GoBack", $notePara.Range)
